# Auto-generated: update cached market-price / profit values per commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 840.7273
$ws.Range("I11").Value = 840.7273
$ws.Range("K11").Value = 840.7273
$ws.Range("M11").Value = -700.7273
$ws.Range("H49").Value = 1499
$ws.Range("J49").Value = 1499
$ws.Range("L49").Value = 4497
$ws.Range("N49").Value = -4769
$ws.Range("H64").Value = 8270
$ws.Range("I64").Value = 3720
$ws.Range("J64").Value = 11682.5
$ws.Range("K64").Value = 3720
$ws.Range("L64").Value = 11682.5
$ws.Range("M64").Value = -3472
$ws.Range("N64").Value = -12178.5
$ws.Range("H67").Value = 8270
$ws.Range("I67").Value = 3720
$ws.Range("J67").Value = 11682.5
$ws.Range("K67").Value = 3720
$ws.Range("L67").Value = 11682.5
$ws.Range("M67").Value = -2862
$ws.Range("N67").Value = -13398.5
$ws.Range("H106").Value = 4233.2144
$ws.Range("I106").Value = 3672.75
$ws.Range("K106").Value = 3672.75
$ws.Range("M106").Value = -3041.75
$ws.Range("H137").Value = 4096.9033
$ws.Range("I137").Value = 4346.731
$ws.Range("J137").Value = 2797.8
$ws.Range("K137").Value = 13040.193
$ws.Range("L137").Value = 8393.400000000001
$ws.Range("M137").Value = -10490.193
$ws.Range("N137").Value = -13493.4
$ws.Range("H138").Value = 2405.152
$ws.Range("J138").Value = 3546.875
$ws.Range("L138").Value = 10640.625
$ws.Range("N138").Value = -20920.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7483.3335
$ws.Range("I32").Value = 5765.857
$ws.Range("K32").Value = 5765.857
$ws.Range("M32").Value = -5478.857
$ws.Range("H61").Value = 5963.4634
$ws.Range("I61").Value = 6084.846
$ws.Range("K61").Value = 6084.846
$ws.Range("M61").Value = -5872.846
$ws.Range("H74").Value = 4247.591
$ws.Range("I74").Value = 3265.6316
$ws.Range("J74").Value = 10466.667
$ws.Range("K74").Value = 3265.6316
$ws.Range("L74").Value = 10466.667
$ws.Range("M74").Value = -2391.6316
$ws.Range("N74").Value = -12214.667
$ws.Range("H77").Value = 4247.591
$ws.Range("I77").Value = 3265.6316
$ws.Range("J77").Value = 10466.667
$ws.Range("K77").Value = 16328.158
$ws.Range("L77").Value = 52333.335
$ws.Range("M77").Value = -11960.158
$ws.Range("N77").Value = -61069.335
$ws.Range("H102").Value = 4191
$ws.Range("I102").Value = 2889.5
$ws.Range("J102").Value = 12000
$ws.Range("K102").Value = 2889.5
$ws.Range("L102").Value = 12000
$ws.Range("M102").Value = -1267.5
$ws.Range("N102").Value = -15244
$ws.Range("H136").Value = 5963.4634
$ws.Range("I136").Value = 6084.846
$ws.Range("K136").Value = 18254.538
$ws.Range("M136").Value = -15704.538

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1737.25
$ws.Range("I22").Value = 1983
$ws.Range("K22").Value = 1983
$ws.Range("M22").Value = -1810
$ws.Range("H86").Value = 2168.85
$ws.Range("I86").Value = 1885.9333
$ws.Range("K86").Value = 1885.9333
$ws.Range("M86").Value = -762.9332999999999
$ws.Range("H89").Value = 2168.85
$ws.Range("I89").Value = 1885.9333
$ws.Range("K89").Value = 9429.666499999999
$ws.Range("M89").Value = -3813.666499999999
$ws.Range("H94").Value = 1202.1765
$ws.Range("I94").Value = 434.7143
$ws.Range("K94").Value = 434.7143
$ws.Range("M94").Value = 16.28570000000002
$ws.Range("H105").Value = 1478
$ws.Range("I105").Value = 1475.9375
$ws.Range("K105").Value = 1475.9375
$ws.Range("M105").Value = 271.0625
$ws.Range("H134").Value = 5866.735
$ws.Range("I134").Value = 5665.7383
$ws.Range("K134").Value = 16997.2149
$ws.Range("M134").Value = -14462.2149

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 12869.941
$ws.Range("I58").Value = 9594
$ws.Range("K58").Value = 9594
$ws.Range("M58").Value = -9391
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H100").Value = 72999.336
$ws.Range("J100").Value = 84999.5
$ws.Range("L100").Value = 84999.5
$ws.Range("N100").Value = -87163.5
$ws.Range("H105").Value = 55556812
$ws.Range("I105").Value = 83334536
$ws.Range("K105").Value = 83334536
$ws.Range("M105").Value = -83332789
$ws.Range("H122").Value = 4041.15
$ws.Range("I122").Value = 4115.3076
$ws.Range("K122").Value = 12345.9228
$ws.Range("M122").Value = -9895.9228
$ws.Range("H136").Value = 12869.941
$ws.Range("I136").Value = 9594
$ws.Range("K136").Value = 28782
$ws.Range("M136").Value = -26232

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 21221988
$ws.Range("I4").Value = 17689486
$ws.Range("K4").Value = 53068458
$ws.Range("M4").Value = -53068346
$ws.Range("H12").Value = 241.6
$ws.Range("I12").Value = 11
$ws.Range("K12").Value = 33
$ws.Range("M12").Value = 140
$ws.Range("H22").Value = 1389.6552
$ws.Range("J22").Value = 1496.5385
$ws.Range("L22").Value = 4489.6155
$ws.Range("N22").Value = -4827.6155
$ws.Range("H26").Value = 114158.11
$ws.Range("J26").Value = 3906.6428
$ws.Range("L26").Value = 11719.9284
$ws.Range("N26").Value = -12295.9284
$ws.Range("H27").Value = 1389.6552
$ws.Range("J27").Value = 1496.5385
$ws.Range("L27").Value = 4489.6155
$ws.Range("N27").Value = -4693.6155
$ws.Range("H132").Value = 2111.1538
$ws.Range("J132").Value = 2141.6667
$ws.Range("L132").Value = 19275.0003
$ws.Range("N132").Value = -24335.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 738.375
$ws.Range("I97").Value = 690.087
$ws.Range("J97").Value = 861.7778
$ws.Range("K97").Value = 690.087
$ws.Range("L97").Value = 861.7778
$ws.Range("M97").Value = -194.087
$ws.Range("N97").Value = -1853.7778
$ws.Range("H132").Value = 5867.2
$ws.Range("I132").Value = 5120.278
$ws.Range("J132").Value = 7787.857
$ws.Range("K132").Value = 15360.834
$ws.Range("L132").Value = 23363.571
$ws.Range("M132").Value = -12830.834
$ws.Range("N132").Value = -28423.571

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7997
$ws.Range("I7").Value = 7997.8
$ws.Range("J7").Value = 7993
$ws.Range("K7").Value = 7997.8
$ws.Range("L7").Value = 7993
$ws.Range("M7").Value = -7885.8
$ws.Range("N7").Value = -8217
$ws.Range("H40").Value = 4861.375
$ws.Range("I40").Value = 4797.3335
$ws.Range("K40").Value = 4797.3335
$ws.Range("M40").Value = -4661.3335
$ws.Range("H93").Value = 19819.785
$ws.Range("I93").Value = 8498.143
$ws.Range("J93").Value = 31141.428
$ws.Range("K93").Value = 8498.143
$ws.Range("L93").Value = 31141.428
$ws.Range("M93").Value = -7250.143
$ws.Range("N93").Value = -33637.428
$ws.Range("H100").Value = 4461.8823
$ws.Range("I100").Value = 2606.5
$ws.Range("K100").Value = 2606.5
$ws.Range("M100").Value = -2065.5
$ws.Range("H122").Value = 2864.1667
$ws.Range("I122").Value = 2877.4
$ws.Range("J122").Value = 2798
$ws.Range("K122").Value = 8632.200000000001
$ws.Range("L122").Value = 8394
$ws.Range("M122").Value = -6182.200000000001
$ws.Range("N122").Value = -13294
$ws.Range("H126").Value = 7997
$ws.Range("I126").Value = 7997.8
$ws.Range("J126").Value = 7993
$ws.Range("K126").Value = 23993.4
$ws.Range("L126").Value = 23979
$ws.Range("M126").Value = -21523.4
$ws.Range("N126").Value = -28919
$ws.Range("H136").Value = 90917160
$ws.Range("I136").Value = 62508440
$ws.Range("K136").Value = 187525320
$ws.Range("M136").Value = -187522770

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1969.25
$ws.Range("I96").Value = 938.5
$ws.Range("K96").Value = 938.5
$ws.Range("M96").Value = 434.5
$ws.Range("H100").Value = 1550
$ws.Range("I100").Value = 1575
$ws.Range("J100").Value = 1500
$ws.Range("K100").Value = 3150
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -2609
$ws.Range("N100").Value = -4082
$ws.Range("H122").Value = 2867
$ws.Range("I122").Value = 1739.8182
$ws.Range("K122").Value = 5219.4546
$ws.Range("M122").Value = -2769.4546
$ws.Range("H132").Value = 7219.625
$ws.Range("I132").Value = 6251
$ws.Range("J132").Value = 14000
$ws.Range("K132").Value = 18753
$ws.Range("L132").Value = 42000
$ws.Range("M132").Value = -16223
$ws.Range("N132").Value = -47060
$ws.Range("H136").Value = 10179.385
$ws.Range("I136").Value = 9303
$ws.Range("K136").Value = 27909
$ws.Range("M136").Value = -25359
